# HTML REPORT GENERATOR - PART3
# Appends the latest "145_beta" sprint run to the BETA sheet's history table
# and refreshes the previous run's captured timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BETA")

# The existing row's "Run Time" (B2) gets re-stamped with a (sub-millisecond)
# refined timestamp value.
$ws.Range("B2").Value = 44355.68716612268

# New row 3: the 2021-06-09 "145_beta" sprint run.
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2021-06-09"
$ws.Range("A3").Style = "Normal"

$ws.Range("B3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B3").Value = 44356.6332471106

$ws.Range("C3").Value = "145_beta"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "105"
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = 105
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 2.857780133333333
